# Update cryptos list with new price/volume snapshot values
# Two swapped-row reorderings (NEARProtocol<->ImmutableX, FLOKI<->Mantle) are
# included as straightforward cell text overwrites, matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.972.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.19%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.798.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.09%  "

# Row 4
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "695.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +10.68%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.29%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.797.81"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.12%  "

# Row 8
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.53%  "

# Row 10
$ws.Range("E10").Value = "  +3.45%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.49"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +8.78%  "

# Row 12
$ws.Range("E12").Value = "  +1.46%  "

# Row 13
$ws.Range("E13").Value = "  +9.76%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.62%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.441.81"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.23%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.801.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.06%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "70.984.01"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.31%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.88"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.74%  "

# Row 19
$ws.Range("E19").Value = "  +3.30%  "

# Row 20
$ws.Range("E20").Value = "  +0.64%  "

# Row 21
$ws.Range("E21").Value = "  +17.26%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "484.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.74%  "

# Row 23
$ws.Range("E23").Value = "  +2.03%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.37%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000145"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.28%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.05%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.36%  "

# Row 28
$ws.Range("E28").Value = "  +3.18%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.950.73"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.15%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"

# Row 31
$ws.Range("E31").Value = "  +14.39%  "

# Row 32
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.42%  "

# Row 33
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.52"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.00%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.72%  "

# Row 35
$ws.Range("E35").Value = "  +1.75%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.24"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.53%  "

# Row 37
$ws.Range("E37").Value = "  +0.29%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.750.30"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.09%  "

# Row 39
$ws.Range("E39").Value = "  +2.76%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.51"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.85%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.98"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.20%  "

# Row 42
$ws.Range("E42").Value = "  +13.99%  "

# Row 43
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.972"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.92%  "

# Row 44
$ws.Range("B44").Value = "FLOKI"
$ws.Range("C44").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.000326"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +23.09%  "

# Row 45
$ws.Range("E45").Value = "  +0.10%  "

# Row 46
$ws.Range("E46").Value = "  +0.01%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "163.57"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.87%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "49.35"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.11%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "44.88"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.07%  "

# Row 50
$ws.Range("E50").Value = "  +3.16%  "

# Row 51
$ws.Range("E51").Value = "  -1.12%  "
